# Update cryptocurrency price/volume data per scheduled GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" "20.562.53"
Set-TextCell $ws "E2" "  +1.65%  "

Set-TextCell $ws "D3" "1.470.18"
Set-TextCell $ws "E3" "  +1.95%  "

Set-TextCell $ws "E4" "  +0.32%  "

Set-TextCell $ws "D5" "0.9598"
Set-TextCell $ws "E5" "  +4.36%  "

Set-TextCell $ws "D6" "276.91"
Set-TextCell $ws "E6" "  +0.87%  "

Set-TextCell $ws "D7" "0.3565"
Set-TextCell $ws "E7" "  -1.78%  "

Set-TextCell $ws "D8" "0.3062"
Set-TextCell $ws "E8" "  -0.38%  "

Set-TextCell $ws "D9" "1.085"
Set-TextCell $ws "E9" "  +6.58%  "

Set-TextCell $ws "D10" "39.48"
Set-TextCell $ws "E10" "  +1.81%  "

Set-TextCell $ws "D11" "0.06621"
Set-TextCell $ws "E11" "  +2.18%  "

Set-TextCell $ws "E12" "  +0.51%  "

Set-TextCell $ws "D13" "5.452"
Set-TextCell $ws "E13" "  +2.44%  "

Set-TextCell $ws "D14" "18.08"
Set-TextCell $ws "E14" "  +3.95%  "

Set-TextCell $ws "D15" "6.169"
Set-TextCell $ws "E15" "  +2.50%  "

Set-TextCell $ws "D16" "0.9597"
Set-TextCell $ws "E16" "  +2.53%  "

Set-TextCell $ws "D17" "0.00001018"
Set-TextCell $ws "E17" "  +1.10%  "

Set-TextCell $ws "D18" "1.471.18"
Set-TextCell $ws "E18" "  +2.14%  "

Set-TextCell $ws "D19" "0.05953"
Set-TextCell $ws "E19" "  +6.21%  "

Set-TextCell $ws "D20" "68.98"
Set-TextCell $ws "E20" "  +2.31%  "

Set-TextCell $ws "D21" "5.478"
Set-TextCell $ws "E21" "  +2.02%  "

Set-TextCell $ws "E22" "  +2.06%  "

Set-TextCell $ws "D23" "11.22"
Set-TextCell $ws "E23" "  +4.16%  "

Set-TextCell $ws "D24" "2.278"
Set-TextCell $ws "E24" "  +1.66%  "

Set-TextCell $ws "D25" "20.575.44"
Set-TextCell $ws "E25" "  +1.61%  "

Set-TextCell $ws "D26" "145.05"
Set-TextCell $ws "E26" "  +4.75%  "

Set-TextCell $ws "D27" "2.086"
Set-TextCell $ws "E27" "  +1.51%  "

Set-TextCell $ws "D28" "17.09"
Set-TextCell $ws "E28" "  +1.35%  "

Set-TextCell $ws "D29" "1.631.08"
Set-TextCell $ws "E29" "  +2.41%  "

Set-TextCell $ws "D30" "114.70"
Set-TextCell $ws "E30" "  +4.34%  "

Set-TextCell $ws "D31" "3.857"
Set-TextCell $ws "E31" "  -3.19%  "

Set-TextCell $ws "D32" "4.904"
Set-TextCell $ws "E32" "  +1.76%  "

Set-TextCell $ws "D33" "0.07940"
Set-TextCell $ws "E33" "  +4.08%  "

Set-TextCell $ws "D34" "0.7921"
Set-TextCell $ws "E34" "  +0.87%  "

Set-TextCell $ws "D35" "1.231"
Set-TextCell $ws "E35" "  +8.31%  "

Set-TextCell $ws "E36" "  -1.13%  "

Set-TextCell $ws "D37" "0.05708"
Set-TextCell $ws "E37" "  -1.34%  "

Set-TextCell $ws "D38" "4.693"
Set-TextCell $ws "E38" "  +1.27%  "

Set-TextCell $ws "D39" "0.9601"
Set-TextCell $ws "E39" "  +3.03%  "

Set-TextCell $ws "D40" "0.02028"
Set-TextCell $ws "E40" "  +2.45%  "

Set-TextCell $ws "D41" "10.25"
Set-TextCell $ws "E41" "  +1.42%  "

Set-TextCell $ws "D42" "0.1849"
Set-TextCell $ws "E42" "  +0.45%  "

Set-TextCell $ws "D43" "7.287"
Set-TextCell $ws "E43" "  +4.20%  "

Set-TextCell $ws "D44" "3.511"
Set-TextCell $ws "E44" "  +1.12%  "

Set-TextCell $ws "D45" "0.5214"
Set-TextCell $ws "E45" "  +0.72%  "

Set-TextCell $ws "D46" "12.08"
Set-TextCell $ws "E46" "  +3.25%  "

Set-TextCell $ws "D47" "119.36"
Set-TextCell $ws "E47" "  +3.15%  "

Set-TextCell $ws "D48" "0.5163"
Set-TextCell $ws "E48" "  +1.64%  "

Set-TextCell $ws "D49" "1.798"
Set-TextCell $ws "E49" "  +4.23%  "

Set-TextCell $ws "D50" "0.06429"
Set-TextCell $ws "E50" "  +1.20%  "

Set-TextCell $ws "D51" "0.9902"
Set-TextCell $ws "E51" "  +0.22%  "
